$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HWCheck")
$ws.Activate()
Write-Output "Top=$($excel.ActiveWindow.Top)"
Write-Output "Width=$($excel.ActiveWindow.Width)"
Write-Output "Height=$($excel.ActiveWindow.Height)"
